# Update Fitness values (column C) for rows 2-128 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C3").Value = 12884
$ws.Range("C4:C7").Value = 11265
$ws.Range("C8:C10").Value = 9759
$ws.Range("C11:C12").Value = 9170
$ws.Range("C13:C16").Value = 8235
$ws.Range("C17:C21").Value = 8202
$ws.Range("C22:C26").Value = 8070
$ws.Range("C27:C57").Value = 7917
$ws.Range("C58:C71").Value = 7884
$ws.Range("C72:C78").Value = 7867
$ws.Range("C79:C128").Value = 7622
